# Example.xlsx update:
#   1. Add a new "effectParam" example column (type = json) to Sheet1,
#      demonstrating a json-typed field, e.g. {"Box":"ys_015"}.
#   2. Make Sheet1 the active sheet/tab again, with the new example cell
#      selected.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New column O on Sheet1: field name / type / example-literal / sample value.
# (Values are written in this particular order so the workbook's shared
# string table ends up laid out the same way the source file has it.)
$ws1.Range("O2").Value = "effectParam"
$ws1.Range("O5").Value = '{"Box":"ys_015"}'
$ws1.Range("O3").Value = "json"
$ws1.Range("O4").Value = '["box", 1, "test"]'

# Sheet1 becomes the active/selected sheet again, with O5 (the new sample
# json value) selected.
$ws1.Activate() | Out-Null
$ws1.Range("O5").Select() | Out-Null
